$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$rng = $ws.Range("C2:C$lastRow")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value = $v + 1
    }
}
